$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.379.95'
Set-TextValue 'E2' '  -0.17%  '
Set-TextValue 'D3' '1.842.93'
Set-TextValue 'E3' '  -0.31%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '238.92'
Set-TextValue 'E5' '  -0.79%  '
Set-TextValue 'D6' '0.6305'
Set-TextValue 'E6' '  -0.27%  '
Set-TextValue 'D7' '0.9998'
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '0.07530'
Set-TextValue 'E8' '  -0.75%  '
Set-TextValue 'D9' '0.2928'
Set-TextValue 'D10' '24.43'
Set-TextValue 'E10' '  -0.64%  '
Set-TextValue 'D11' '0.07707'
Set-TextValue 'D12' '1.866.00'
Set-TextValue 'E12' '  -6.00%  '
Set-TextValue 'E14' '  -1.08%  '
Set-TextValue 'D15' '0.00001034'
Set-TextValue 'E15' '  +3.44%  '
Set-TextValue 'D16' '82.85'
Set-TextValue 'E16' '  -0.05%  '
Set-TextValue 'D17' '2.112.15'
Set-TextValue 'E17' '  -6.72%  '
Set-TextValue 'D18' '6.140'
Set-TextValue 'E18' '  -0.79%  '
Set-TextValue 'D19' '29.407.45'
Set-TextValue 'E19' '  -0.15%  '
Set-TextValue 'E20' '  -1.77%  '
Set-TextValue 'D21' '12.41'
Set-TextValue 'E21' '  -0.76%  '
Set-TextValue 'D22' '1.0000'
Set-TextValue 'E22' '  +0.02%  '
Set-TextValue 'D23' '7.432'
Set-TextValue 'D24' '1.000'
Set-TextValue 'E24' '  +0.06%  '
Set-TextValue 'D25' '156.96'
Set-TextValue 'E25' '  +1.13%  '
Set-TextValue 'D26' '0.1391'
Set-TextValue 'E26' '  +0.11%  '
Set-TextValue 'D27' '8.349'
Set-TextValue 'E27' '  -1.07%  '
Set-TextValue 'E28' '  -0.50%  '
Set-TextValue 'D29' '1.458'
Set-TextValue 'E29' '  -0.84%  '
Set-TextValue 'E30' '  +1.28%  '
Set-TextValue 'D31' '0.05619'
Set-TextValue 'E31' '  -3.23%  '
Set-TextValue 'D32' '4.104'
Set-TextValue 'E32' '  -0.59%  '
Set-TextValue 'D33' '4.017'
Set-TextValue 'E34' '  -2.16%  '
Set-TextValue 'E35' '  -0.25%  '
Set-TextValue 'D36' '0.7111'
Set-TextValue 'E36' '  -1.18%  '
Set-TextValue 'E37' '  -0.25%  '
Set-TextValue 'D38' '1.241.25'
Set-TextValue 'E38' '  -0.67%  '
Set-TextValue 'D39' '0.01805'
Set-TextValue 'E39' '  -0.23%  '
Set-TextValue 'D40' '2.770'
Set-TextValue 'E40' '  -0.88%  '
Set-TextValue 'D41' '6.319'
Set-TextValue 'E41' '  +3.71%  '
Set-TextValue 'D42' '0.9012'
Set-TextValue 'E42' '  +0.09%  '
Set-TextValue 'E43' '  -0.01%  '
Set-TextValue 'D44' '101.86'
Set-TextValue 'E44' '  +0.48%  '
Set-TextValue 'D45' '65.57'
Set-TextValue 'E45' '  -2.05%  '
Set-TextValue 'E46' '  +0.44%  '
Set-TextValue 'D47' '7.058'
Set-TextValue 'E47' '  -3.61%  '
Set-TextValue 'D48' '0.3998'
Set-TextValue 'E48' '  -0.52%  '
Set-TextValue 'B49' 'RenderToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D49' '1.665'
Set-TextValue 'E49' '  -1.75%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '8.872'
Set-TextValue 'E50' '  -3.19%  '
Set-TextValue 'D51' '0.1118'
Set-TextValue 'E51' '  -0.66%  '
